$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Chapter 19 (JSON) and Chapter 20 (Error handling) rows as DONE,
# matching the formatting used by the rows above (copy style from B19).
$ws.Range("B19").Copy() | Out-Null
$ws.Range("B20:B21").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B20").Value = "DONE"
$ws.Range("B21").Value = "DONE"

# Update the active selection to match the author's final view.
$ws.Range("F17").Select() | Out-Null
